{"js": "// Remove the post-condition hyperlink text \"La postulaci\u00f3n queda registrada\n// en el sistema.\" (a TOC-style entry linking to _bookmark7) and its page\n// number \"5\" from the table cells that still host it, leaving the empty\n// paragraph (and its paragraph mark formatting) in place \u2014 matching the\n// author's prototype cleanup described in the commit message.\n\nconst targetText = \"La postulaci\u00f3n queda registrada en el sistema.\";\n\nconst hits = context.document.body.search(targetText, { matchCase: true });\nhits.load(\"text\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  const hit = hits.items[0];\n\n  // Locate the sibling cell (same row) that carries the page-number run (\"5\").\n  const cell = hit.parentTableCell;\n  const row = cell.parentRow;\n  const lastCell = row.cells.getLast();\n\n  const pageNumberHits = lastCell.body.search(\"5\", { matchCase: true });\n  pageNumberHits.load(\"text\");\n  await context.sync();\n\n  // Clear the hyperlinked text, keeping the paragraph mark intact.\n  hit.insertText(\"\", Word.InsertLocation.replace);\n\n  // Clear the corresponding page-number run in the neighboring cell.\n  if (pageNumberHits.items.length > 0) {\n    pageNumberHits.items[0].insertText(\"\", Word.InsertLocation.replace);\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the post-condition hyperlink text \"La postulaci\u00f3n queda registrada\n# en el sistema.\" (a TOC-style entry linking to _bookmark7) and its page\n# number \"5\" from the table cells that still host it, leaving the empty\n# paragraph (and its paragraph mark formatting) in place \u2014 matching the\n# author's prototype cleanup described in the commit message.\n\n$d = $word.ActiveDocument\n\n$targetText = \"La postulaci\u00f3n queda registrada en el sistema.\"\n\n$findRange = $d.Content\n$find = $findRange.Find\n$find.Text = $targetText\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif ($found -and $find.Found) {\n    # Locate the table cell/row that holds the matched hyperlink text so we\n    # can also clear the neighboring page-number cell (\"5\").\n    $cell = $findRange.Cells(1)\n    $table = $cell.Tables(1)\n    $row = $table.Rows($cell.RowIndex)\n    $lastCell = $row.Cells($row.Cells.Count)\n\n    # Clear the page-number run (\"5\") in the sibling cell first (deleting the\n    # found range first would shift the document and is avoided here by\n    # resolving the sibling cell up front).\n    $lastCell.Range.Delete()\n\n    # Remove the hyperlinked text itself, leaving the paragraph mark intact.\n    $findRange.Delete()\n}\n"}
